$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B4 was stored as a text "2"; change it to a genuine numeric value 2.
$ws.Range("B4").Value = 2

# Add the new annotation row (row 5).
$ws.Range("A5").Value = "Sunsi Wu"

# B5 must stay a text string "3" (not be auto-coerced into a number like B4).
# Temporarily force a text number format so the "3" is kept as a string, then
# clear the formatting again so the cell ends up unstyled (matching the
# other data cells) while remaining text-typed.
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "3"
$ws.Range("B5").ClearFormats()

$ws.Range("C5").Value = "what about"
$ws.Range("D5").Value = "QSN"
$ws.Range("E5").Value = "MET"
$ws.Range("F5").Value = "53dcf950-aee9-43ba-bb93-9e7c5cd5833d"
$ws.Range("G5").Value = "By5SY2gA-_annotated.xlsx"
$ws.Range("H5").Value = "For instance, what about averaging WordNet path-based distance metrics and distance in word embedding space (for word similarity), and other ways of applying the affect data to email tone prediction?"
